$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.030.77"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "2.103.02"
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "'349.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.68%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "'0.5164"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("D8").Value = "'0.4439"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "'52.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.49%  "
$ws.Range("D10").Value = "'0.08978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "'1.171"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("E12").Value = "  +4.19%  "
$ws.Range("D13").Value = "2.102.02"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "'8.228"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'6.735"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Value = "'98.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "'0.00001148"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'20.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.74%  "
$ws.Range("D20").Value = "'0.06675"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "'6.226"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").Value = "30.144.53"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").Value = "'2.343"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "2.351.61"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "'21.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").Value = "'2.548"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'162.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "'133.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "'1.176"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("D32").Value = "'0.1062"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").Value = "'1.647"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'6.220"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "'3.968"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "'5.922"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'10.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").Value = "'0.02572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'0.2298"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").Value = "'12.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.325"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.55%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6804"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "'14.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.43%  "
$ws.Range("D45").Value = "'0.6369"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").Value = "'2.289"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'3.642"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "'1.218"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").Value = "'82.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "'0.07228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
